{"js": "// Appends six new FAQ paragraphs (three bold question / answer pairs about\n// \"corte de control\" and \"modularizar\") to the end of the document, matching\n// the existing \"Prrafodelista\" / numId=1 bullet list used throughout the doc\n// (ilvl 0 = bold question, ilvl 1 = plain answer).\n//\n// Word's Office.js Range/Body API has no way to set the paragraph-mark run\n// properties (the <w:rPr> that lives directly under <w:pPr>, as opposed to\n// the run's own <w:rPr>) that Word itself stamps onto a bold heading\n// paragraph, so we build the exact WordprocessingML for the new paragraphs\n// (the same fragment a real Word session would have produced) and hand it to\n// insertOoxml, which is the Office.js entry point for Range.InsertXML.\nconst newParagraphsXml = \"<?xml version=\\\"1.0\\\" encoding=\\\"UTF-8\\\" standalone=\\\"yes\\\"?><pkg:package xmlns:pkg=\\\"http://schemas.microsoft.com/office/2006/xmlPackage\\\"><pkg:part pkg:name=\\\"/word/document.xml\\\" pkg:contentType=\\\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\\\"><pkg:xmlData><w:document xmlns:w=\\\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\\\"><w:body><w:p><w:pPr><w:pStyle w:val=\\\"Prrafodelista\\\"/><w:numPr><w:ilvl w:val=\\\"0\\\"/><w:numId w:val=\\\"1\\\"/></w:numPr><w:rPr><w:b/><w:lang w:val=\\\"es-AR\\\"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:lang w:val=\\\"es-AR\\\"/></w:rPr><w:t>\\u00bfCu\\u00e1ndo \\u201cuso corte de control\\u201d?</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\\\"Prrafodelista\\\"/><w:numPr><w:ilvl w:val=\\\"1\\\"/><w:numId w:val=\\\"1\\\"/></w:numPr><w:rPr><w:lang w:val=\\\"es-AR\\\"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=\\\"es-AR\\\"/></w:rPr><w:t>Generalmente los algoritmos de corte de control se pueden usar cuando los datos vienen ordenados bajo cierto criterio.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\\\"Prrafodelista\\\"/><w:numPr><w:ilvl w:val=\\\"0\\\"/><w:numId w:val=\\\"1\\\"/></w:numPr><w:rPr><w:b/><w:lang w:val=\\\"es-AR\\\"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:lang w:val=\\\"es-AR\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\">\\u00bfCu\\u00e1ndo </w:t></w:r><w:proofErr w:type=\\\"spellStart\\\"/><w:r><w:rPr><w:b/><w:lang w:val=\\\"es-AR\\\"/></w:rPr><w:t>modularizo</w:t></w:r><w:proofErr w:type=\\\"spellEnd\\\"/><w:r><w:rPr><w:b/><w:lang w:val=\\\"es-AR\\\"/></w:rPr><w:t>?</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\\\"Prrafodelista\\\"/><w:numPr><w:ilvl w:val=\\\"1\\\"/><w:numId w:val=\\\"1\\\"/></w:numPr><w:rPr><w:lang w:val=\\\"es-AR\\\"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=\\\"es-AR\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\">Hay que pensar que un m\\u00f3dulo (ya sea funci\\u00f3n o procedimiento) debe ser lo m\\u00e1s gen\\u00e9rico y reutilizable posible. Lo ideal es que solo se encargue de hacer una tarea en </w:t></w:r><w:proofErr w:type=\\\"gramStart\\\"/><w:r><w:rPr><w:lang w:val=\\\"es-AR\\\"/></w:rPr><w:t>especifico</w:t></w:r><w:proofErr w:type=\\\"gramEnd\\\"/><w:r><w:rPr><w:lang w:val=\\\"es-AR\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\">. Si el m\\u00f3dulo que hiciste hace </w:t></w:r><w:proofErr w:type=\\\"spellStart\\\"/><w:proofErr w:type=\\\"gramStart\\\"/><w:r><w:rPr><w:lang w:val=\\\"es-AR\\\"/></w:rPr><w:t>mas</w:t></w:r><w:proofErr w:type=\\\"spellEnd\\\"/><w:proofErr w:type=\\\"gramEnd\\\"/><w:r><w:rPr><w:lang w:val=\\\"es-AR\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\"> de una tarea, entonces se puede seguir </w:t></w:r><w:proofErr w:type=\\\"spellStart\\\"/><w:r><w:rPr><w:lang w:val=\\\"es-AR\\\"/></w:rPr><w:t>modularizando</w:t></w:r><w:proofErr w:type=\\\"spellEnd\\\"/><w:r><w:rPr><w:lang w:val=\\\"es-AR\\\"/></w:rPr><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\\\"Prrafodelista\\\"/><w:numPr><w:ilvl w:val=\\\"0\\\"/><w:numId w:val=\\\"1\\\"/></w:numPr><w:rPr><w:b/><w:lang w:val=\\\"es-AR\\\"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:lang w:val=\\\"es-AR\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\">\\u00bfVale la pena </w:t></w:r><w:proofErr w:type=\\\"spellStart\\\"/><w:r><w:rPr><w:b/><w:lang w:val=\\\"es-AR\\\"/></w:rPr><w:t>modularizar</w:t></w:r><w:proofErr w:type=\\\"spellEnd\\\"/><w:r><w:rPr><w:b/><w:lang w:val=\\\"es-AR\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\"> una suma?</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\\\"Prrafodelista\\\"/><w:numPr><w:ilvl w:val=\\\"1\\\"/><w:numId w:val=\\\"1\\\"/></w:numPr><w:rPr><w:lang w:val=\\\"es-AR\\\"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=\\\"es-AR\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\">Si esas peque\\u00f1as l\\u00edneas de c\\u00f3digo las usas muchas veces a lo largo de tu programa, entonces tal vez podr\\u00edas </w:t></w:r><w:proofErr w:type=\\\"spellStart\\\"/><w:r><w:rPr><w:lang w:val=\\\"es-AR\\\"/></w:rPr><w:t>modularizarlas</w:t></w:r><w:proofErr w:type=\\\"spellEnd\\\"/><w:r><w:rPr><w:lang w:val=\\\"es-AR\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\"> para que quede </w:t></w:r><w:proofErr w:type=\\\"spellStart\\\"/><w:r><w:rPr><w:lang w:val=\\\"es-AR\\\"/></w:rPr><w:t>mas</w:t></w:r><w:proofErr w:type=\\\"spellEnd\\\"/><w:r><w:rPr><w:lang w:val=\\\"es-AR\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\"> legible el c\\u00f3digo final.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\";\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Append a fresh, empty paragraph strictly AFTER the current last paragraph\n// first. Inserting the OOXML directly at the body/paragraph \"End\" caret\n// would land *inside* the last paragraph (right before its own paragraph\n// mark) and merge with it, destroying the existing final bullet. Using a\n// brand-new placeholder paragraph as the insertion target keeps the\n// original last paragraph completely untouched.\nconst placeholder = lastParagraph.insertParagraph(\"\", \"After\");\nawait context.sync();\n\n// Replace that placeholder paragraph's content with the new OOXML fragment.\nconst placeholderRange = placeholder.getRange(\"Whole\");\nplaceholderRange.insertOoxml(newParagraphsXml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# WordprocessingML fragment for the new FAQ paragraphs (questions in bold at\n# ilvl 0, answers at ilvl 1), wrapped as a Flat OPC package so that\n# Range.InsertXML accepts it and reproduces the exact OOXML the diff adds.\n$newParagraphs = @'\n<w:p><w:pPr><w:pStyle w:val=\"Prrafodelista\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr><w:rPr><w:b/><w:lang w:val=\"es-AR\"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:lang w:val=\"es-AR\"/></w:rPr><w:t>&#191;Cu&#225;ndo &#8220;uso corte de control&#8221;?</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Prrafodelista\"/><w:numPr><w:ilvl w:val=\"1\"/><w:numId w:val=\"1\"/></w:numPr><w:rPr><w:lang w:val=\"es-AR\"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=\"es-AR\"/></w:rPr><w:t>Generalmente los algoritmos de corte de control se pueden usar cuando los datos vienen ordenados bajo cierto criterio.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Prrafodelista\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr><w:rPr><w:b/><w:lang w:val=\"es-AR\"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:lang w:val=\"es-AR\"/></w:rPr><w:t xml:space=\"preserve\">&#191;Cu&#225;ndo </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:rPr><w:b/><w:lang w:val=\"es-AR\"/></w:rPr><w:t>modularizo</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:rPr><w:b/><w:lang w:val=\"es-AR\"/></w:rPr><w:t>?</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Prrafodelista\"/><w:numPr><w:ilvl w:val=\"1\"/><w:numId w:val=\"1\"/></w:numPr><w:rPr><w:lang w:val=\"es-AR\"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=\"es-AR\"/></w:rPr><w:t xml:space=\"preserve\">Hay que pensar que un m&#243;dulo (ya sea funci&#243;n o procedimiento) debe ser lo m&#225;s gen&#233;rico y reutilizable posible. Lo ideal es que solo se encargue de hacer una tarea en </w:t></w:r><w:proofErr w:type=\"gramStart\"/><w:r><w:rPr><w:lang w:val=\"es-AR\"/></w:rPr><w:t>especifico</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:rPr><w:lang w:val=\"es-AR\"/></w:rPr><w:t xml:space=\"preserve\">. Si el m&#243;dulo que hiciste hace </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:proofErr w:type=\"gramStart\"/><w:r><w:rPr><w:lang w:val=\"es-AR\"/></w:rPr><w:t>mas</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:proofErr w:type=\"gramEnd\"/><w:r><w:rPr><w:lang w:val=\"es-AR\"/></w:rPr><w:t xml:space=\"preserve\"> de una tarea, entonces se puede seguir </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:rPr><w:lang w:val=\"es-AR\"/></w:rPr><w:t>modularizando</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:rPr><w:lang w:val=\"es-AR\"/></w:rPr><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Prrafodelista\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr><w:rPr><w:b/><w:lang w:val=\"es-AR\"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:lang w:val=\"es-AR\"/></w:rPr><w:t xml:space=\"preserve\">&#191;Vale la pena </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:rPr><w:b/><w:lang w:val=\"es-AR\"/></w:rPr><w:t>modularizar</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:rPr><w:b/><w:lang w:val=\"es-AR\"/></w:rPr><w:t xml:space=\"preserve\"> una suma?</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Prrafodelista\"/><w:numPr><w:ilvl w:val=\"1\"/><w:numId w:val=\"1\"/></w:numPr><w:rPr><w:lang w:val=\"es-AR\"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=\"es-AR\"/></w:rPr><w:t xml:space=\"preserve\">Si esas peque&#241;as l&#237;neas de c&#243;digo las usas muchas veces a lo largo de tu programa, entonces tal vez podr&#237;as </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:rPr><w:lang w:val=\"es-AR\"/></w:rPr><w:t>modularizarlas</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:rPr><w:lang w:val=\"es-AR\"/></w:rPr><w:t xml:space=\"preserve\"> para que quede </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:rPr><w:lang w:val=\"es-AR\"/></w:rPr><w:t>mas</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:rPr><w:lang w:val=\"es-AR\"/></w:rPr><w:t xml:space=\"preserve\"> legible el c&#243;digo final.</w:t></w:r></w:p>\n'@\n\n$flatOpc = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' + $newParagraphs + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n# Insert at a collapsed range sitting right at the very end of the document\n# (past the last paragraph's mark) so InsertXML appends the new paragraphs\n# after the existing content instead of splitting the last paragraph.\n$endPos = $d.Content.End\n$insertionRange = $d.Range($endPos, $endPos)\n$insertionRange.InsertXML($flatOpc)\n\n$d.Save()\n"}
